$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining rows (9, 10, 11) that were left incomplete in the
# original workbook. Each row lists the 4 color bands and the resulting
# nominal value + tolerance.
$ws.Range("B9").Value = "laranja"
$ws.Range("C9").Value = "branco"
$ws.Range("D9").Value = "vermelho"
$ws.Range("E9").Value = "dourado"
$ws.Range("H9").Value = "3k9Ω±5%"

$ws.Range("B10").Value = "cinza"
$ws.Range("C10").Value = "vermelho"
$ws.Range("D10").Value = "marrom"
$ws.Range("E10").Value = "dourado"
$ws.Range("H10").Value = "820Ω±5%"

$ws.Range("B11").Value = "marrom"
$ws.Range("C11").Value = "vermelho"
$ws.Range("D11").Value = "laranja"
$ws.Range("E11").Value = "dourado"
$ws.Range("H11").Value = "12kΩ±5%"

# Reset the view so the top-left cell is the default and the last
# selected cell is H11 (matches the author's final save state).
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H11").Select()
